$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.010381937026978
$ws.Range("B1").Value = 1.308330178260803
$ws.Range("C1").Value = 8.954529762268066
$ws.Range("D1").Value = 2.382446765899658
$ws.Range("E1").Value = 1.278698682785034
